$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update price (D) and 1h volume-change (E) figures for the refreshed crypto snapshot ---
$ws.Range("D2").Value = "37.147.00"
$ws.Range("E2").Value = "  +1.21%  "
$ws.Range("D3").Value = "2.056.66"
$ws.Range("E3").Value = "  -3.40%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'248.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.64%  "
$ws.Range("D6").Value = "'0.655"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.90%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'54.92"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +16.55%  "
$ws.Range("D9").Value = "'61.96"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.28%  "
$ws.Range("E10").Value = "  +0.72%  "
$ws.Range("D11").Value = "'0.0787"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.84%  "
$ws.Range("E12").Value = "  +5.48%  "
$ws.Range("D13").Value = "'15.05"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.67%  "
$ws.Range("D14").Value = "2.354.37"
$ws.Range("E14").Value = "  -3.55%  "
$ws.Range("D15").Value = "'0.818"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.93%  "
$ws.Range("D16").Value = "'5.23"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.95%  "
$ws.Range("D17").Value = "2.051.09"
$ws.Range("E17").Value = "  -3.68%  "
$ws.Range("D18").Value = "37.070.82"
$ws.Range("E18").Value = "  +0.91%  "
$ws.Range("D19").Value = "'72.28"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.85%  "
$ws.Range("D20").Value = "0.0₃0903"
$ws.Range("E20").Value = "  +7.73%  "
$ws.Range("D21").Value = "'14.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.07%  "
$ws.Range("E22").Value = "  +2.13%  "
$ws.Range("D23").Value = "'236.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.88%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("E25").Value = "  -2.52%  "
$ws.Range("D26").Value = "'170.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.94%  "
$ws.Range("D27").Value = "'9.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.33%  "
$ws.Range("D28").Value = "'20.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.43%  "
$ws.Range("E29").Value = "  -2.51%  "
$ws.Range("E30").Value = "  -0.40%  "
$ws.Range("D31").Value = "'4.56"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.11%  "
$ws.Range("E32").Value = "  +11.66%  "
$ws.Range("D33").Value = "'0.0625"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.06%  "
$ws.Range("D34").Value = "'4.31"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.28%  "
$ws.Range("D35").Value = "'0.0885"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -8.06%  "
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("D37").Value = "'2.26"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.11%  "
$ws.Range("D38").Value = "'1.75"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -8.02%  "
$ws.Range("E39").Value = "  -0.39%  "
$ws.Range("E40").Value = "  +22.56%  "
$ws.Range("D41").Value = "'18.26"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +12.86%  "
$ws.Range("D42").Value = "'0.0224"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.98%  "
$ws.Range("D43").Value = "'15.12"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -46.70%  "
$ws.Range("E44").Value = "  -5.44%  "
$ws.Range("D48").Value = "1.296.52"
$ws.Range("E48").Value = "  -4.54%  "
$ws.Range("E49").Value = "  +3.38%  "
$ws.Range("E50").Value = "  +2.90%  "
$ws.Range("D51").Value = "'6.78"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.65%  "

# --- Rows 45-47 rotate coin identities (Aave / FTXToken / HuobiToken) along with their data ---
$ws.Range("B45").Value = "FTXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D45").Value = "'4.38"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +46.98%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'96.11"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.11%  "
$ws.Range("B47").Value = "HuobiToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D47").Value = "'2.84"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.01%  "
